# Apply the "adds journals and footnotes" edit to research_notes.xlsx
#
# The edit inserts two new journal entries (rows) into the sorted table on
# Sheet1, right after "Canadian Journal of Political Science" (h5_index 25)
# and before "International Journal of Public Opinion Research" (h5_index 24):
#   - Australian Journal of Public Administration (h5_index 25)
#   - Acta Politica (h5_index 24)
# All later rows shift down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 38 (pushes the old row 38.."International
# Journal of Public Opinion Research" and everything after it down by two).
$ws.Rows(38).Insert()
$ws.Rows(38).Insert()

# Fill row 39 first (Acta Politica), then row 38 (Australian Journal of
# Public Administration) so that new shared-string entries are created in
# the same order as the target workbook.
$ws.Range("A39").Value2 = "Acta Politica"
$ws.Range("B39").Value2 = "<a href='https://www.palgrave.com/gp/journal/41269/authors/submission'target='_blank'>Research Note</a>"
$ws.Range("C39").Value2 = "*N/A*"
$ws.Range("D39").Value2 = 24

$ws.Range("A38").Value2 = "Australian Journal of Public Administration"
$ws.Range("B38").Value2 = "<a href='https://onlinelibrary.wiley.com/page/journal/14678500/homepage/forauthors.html'target='_blank'>Research Note</a>"
$ws.Range("C38").Value2 = "4k -- 6k words"
$ws.Range("D38").Value2 = 25

# Restore the workbook/sheet view state (best effort): scroll so row 10 is
# near the top and select H41, matching where the author was last working.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H41").Select()
